$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.269101333333333
$ws.Range("H2").Value = 6.807304
$ws.Range("I2").Value = 0.02891211995713196
$ws.Range("J2").Value = 0.02891211995713196
$ws.Range("M2").Value = 218.721583
$ws.Range("N2").Value = 656.164749
$ws.Range("O2").Value = 0.7793342808141792
$ws.Range("P2").Value = 0.7793342808141792
$ws.Range("Q2").Value = 496.3014356140774
$ws.Range("R2").Value = 4466.712920526697
$ws.Range("S2").Value = 0.02253220621360471
$ws.Range("T2").Value = 0.02253220621360471

$ws.Range("G3").Value = 2.269101333333333
$ws.Range("H3").Value = 6.807304
$ws.Range("I3").Value = 0.02891211995713196
$ws.Range("J3").Value = 0.02891211995713196
$ws.Range("O3").Value = 0.164954193449581
$ws.Range("P3").Value = 0.164954193449581
$ws.Range("Q3").Value = 105.0473526380258
$ws.Range("R3").Value = 945.4261737422321
$ws.Range("S3").Value = 0.004769175428446236
$ws.Range("T3").Value = 0.004769175428446236

$ws.Range("G4").Value = 2.269101333333333
$ws.Range("H4").Value = 6.807304
$ws.Range("I4").Value = 0.02891211995713196
$ws.Range("J4").Value = 0.02891211995713196
$ws.Range("M4").Value = 8.51531
$ws.Range("N4").Value = 25.54593
$ws.Range("O4").Value = 0.03034118948727519
$ws.Range("P4").Value = 0.03034118948727519
$ws.Range("Q4").Value = 19.32210127474666
$ws.Range("R4").Value = 173.89891147272
$ws.Range("S4").Value = 0.0008772281100981713
$ws.Range("T4").Value = 0.0008772281100981712

$ws.Range("G5").Value = 2.269101333333333
$ws.Range("H5").Value = 6.807304
$ws.Range("I5").Value = 0.02891211995713196
$ws.Range("J5").Value = 0.02891211995713196
$ws.Range("M5").Value = 7.120231
$ws.Range("N5").Value = 21.360693
$ws.Range("O5").Value = 0.02537033624896462
$ws.Range("P5").Value = 0.02537033624896462
$ws.Range("Q5").Value = 16.15652565574133
$ws.Range("R5").Value = 145.408730901672
$ws.Range("S5").Value = 0.0007335102049828384
$ws.Range("T5").Value = 0.0007335102049828383

$ws.Range("I6").Value = 0.7238963226334669
$ws.Range("J6").Value = 0.7238963226334669
$ws.Range("M6").Value = 218.721583
$ws.Range("N6").Value = 656.164749
$ws.Range("O6").Value = 0.7793342808141792
$ws.Range("P6").Value = 0.7793342808141792
$ws.Range("Q6").Value = 12426.30373322441
$ws.Range("R6").Value = 111836.7335990197
$ws.Range("S6").Value = 0.564157219983582
$ws.Range("T6").Value = 0.564157219983582

$ws.Range("I7").Value = 0.7238963226334669
$ws.Range("J7").Value = 0.7238963226334669
$ws.Range("O7").Value = 0.164954193449581
$ws.Range("P7").Value = 0.164954193449581
$ws.Range("S7").Value = 0.1194097340411212
$ws.Range("T7").Value = 0.1194097340411212

$ws.Range("I8").Value = 0.7238963226334669
$ws.Range("J8").Value = 0.7238963226334669
$ws.Range("M8").Value = 8.51531
$ws.Range("N8").Value = 25.54593
$ws.Range("O8").Value = 0.03034118948727519
$ws.Range("P8").Value = 0.03034118948727519
$ws.Range("Q8").Value = 483.7832050738367
$ws.Range("R8").Value = 4354.04884566453
$ws.Range("S8").Value = 0.02196387549416371
$ws.Range("T8").Value = 0.02196387549416371

$ws.Range("I9").Value = 0.7238963226334669
$ws.Range("J9").Value = 0.7238963226334669
$ws.Range("M9").Value = 7.120231
$ws.Range("N9").Value = 21.360693
$ws.Range("O9").Value = 0.02537033624896462
$ws.Range("P9").Value = 0.02537033624896462
$ws.Range("Q9").Value = 404.5241070549504
$ws.Range("R9").Value = 3640.716963494553
$ws.Range("S9").Value = 0.01836549311460004
$ws.Range("T9").Value = 0.01836549311460003

$ws.Range("G10").Value = 18.57257166666666
$ws.Range("H10").Value = 55.717715
$ws.Range("I10").Value = 0.2366454120188096
$ws.Range("J10").Value = 0.2366454120188096
$ws.Range("M10").Value = 218.721583
$ws.Range("N10").Value = 656.164749
$ws.Range("O10").Value = 0.7793342808141792
$ws.Range("P10").Value = 0.7793342808141792
$ws.Range("Q10").Value = 4062.222275314281
$ws.Range("R10").Value = 36560.00047782854
$ws.Range("S10").Value = 0.1844258819836541
$ws.Range("T10").Value = 0.1844258819836541

$ws.Range("G11").Value = 18.57257166666666
$ws.Range("H11").Value = 55.717715
$ws.Range("I11").Value = 0.2366454120188096
$ws.Range("J11").Value = 0.2366454120188096
$ws.Range("O11").Value = 0.164954193449581
$ws.Range("P11").Value = 0.164954193449581
$ws.Range("Q11").Value = 859.8115282922605
$ws.Range("R11").Value = 7738.303754630345
$ws.Range("S11").Value = 0.03903565307310652
$ws.Range("T11").Value = 0.03903565307310652

$ws.Range("G12").Value = 18.57257166666666
$ws.Range("H12").Value = 55.717715
$ws.Range("I12").Value = 0.2366454120188096
$ws.Range("J12").Value = 0.2366454120188096
$ws.Range("M12").Value = 8.51531
$ws.Range("N12").Value = 25.54593
$ws.Range("O12").Value = 0.03034118948727519
$ws.Range("P12").Value = 0.03034118948727519
$ws.Range("Q12").Value = 158.1512052388833
$ws.Range("R12").Value = 1423.36084714995
$ws.Range("S12").Value = 0.007180103287357011
$ws.Range("T12").Value = 0.007180103287357011

$ws.Range("G13").Value = 18.57257166666666
$ws.Range("H13").Value = 55.717715
$ws.Range("I13").Value = 0.2366454120188096
$ws.Range("J13").Value = 0.2366454120188096
$ws.Range("M13").Value = 7.120231
$ws.Range("N13").Value = 21.360693
$ws.Range("O13").Value = 0.02537033624896462
$ws.Range("P13").Value = 0.02537033624896462
$ws.Range("Q13").Value = 132.2410005307217
$ws.Range("R13").Value = 1190.169004776495
$ws.Range("S13").Value = 0.006003773674691973
$ws.Range("T13").Value = 0.006003773674691973

$ws.Range("G14").Value = 0.8276899999999999
$ws.Range("H14").Value = 2.48307
$ws.Range("I14").Value = 0.01054614539059158
$ws.Range("J14").Value = 0.01054614539059158
$ws.Range("M14").Value = 218.721583
$ws.Range("N14").Value = 656.164749
$ws.Range("O14").Value = 0.7793342808141792
$ws.Range("P14").Value = 0.7793342808141792
$ws.Range("Q14").Value = 181.03366703327
$ws.Range("R14").Value = 1629.30300329943
$ws.Range("S14").Value = 0.008218972633338461
$ws.Range("T14").Value = 0.008218972633338461

$ws.Range("G15").Value = 0.8276899999999999
$ws.Range("H15").Value = 2.48307
$ws.Range("I15").Value = 0.01054614539059158
$ws.Range("J15").Value = 0.01054614539059158
$ws.Range("O15").Value = 0.164954193449581
$ws.Range("P15").Value = 0.164954193449581
$ws.Range("Q15").Value = 38.31765555275666
$ws.Range("R15").Value = 344.85889997481
$ws.Range("S15").Value = 0.001739630906907051
$ws.Range("T15").Value = 0.001739630906907051

$ws.Range("G16").Value = 0.8276899999999999
$ws.Range("H16").Value = 2.48307
$ws.Range("I16").Value = 0.01054614539059158
$ws.Range("J16").Value = 0.01054614539059158
$ws.Range("M16").Value = 8.51531
$ws.Range("N16").Value = 25.54593
$ws.Range("O16").Value = 0.03034118948727519
$ws.Range("P16").Value = 0.03034118948727519
$ws.Range("Q16").Value = 7.048036933899999
$ws.Range("R16").Value = 63.43233240509998
$ws.Range("S16").Value = 0.000319982595656293
$ws.Range("T16").Value = 0.000319982595656293

$ws.Range("G17").Value = 0.8276899999999999
$ws.Range("H17").Value = 2.48307
$ws.Range("I17").Value = 0.01054614539059158
$ws.Range("J17").Value = 0.01054614539059158
$ws.Range("M17").Value = 7.120231
$ws.Range("N17").Value = 21.360693
$ws.Range("O17").Value = 0.02537033624896462
$ws.Range("P17").Value = 0.02537033624896462
$ws.Range("Q17").Value = 5.89334399639
$ws.Range("R17").Value = 53.04009596751
$ws.Range("S17").Value = 0.0002675592546897768
$ws.Range("T17").Value = 0.0002675592546897768
